$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated experiment values for CI (row 2, columns D:AF)
$ws.Range("D2").Value = 1.556312310768631
$ws.Range("E2").Value = 2.198792745987445
$ws.Range("F2").Value = 2.839000178976151
$ws.Range("G2").Value = 3.440493344366906
$ws.Range("H2").Value = 3.984470386134628
$ws.Range("I2").Value = 4.463045771634981
$ws.Range("J2").Value = 4.874540926458894
$ws.Range("K2").Value = 5.220171709781399
$ws.Range("L2").Value = 5.5021343309656
$ws.Range("M2").Value = 5.714910465432629
$ws.Range("N2").Value = 5.861690613148155
$ws.Range("O2").Value = 5.943207533309176
$ws.Range("P2").Value = 5.957422695658019
$ws.Range("Q2").Value = 5.916478857470395
$ws.Range("R2").Value = 5.844445801397892
$ws.Range("S2").Value = 5.757124627193059
$ws.Range("T2").Value = 5.664684741295121
$ws.Range("U2").Value = 5.573471874536271
$ws.Range("V2").Value = 5.487251044832888
$ws.Range("W2").Value = 5.40806261693606
$ws.Range("X2").Value = 5.336812122870487
$ws.Range("Y2").Value = 5.273675972675211
$ws.Range("Z2").Value = 5.218379304801164
$ws.Range("AA2").Value = 5.170384769993985
$ws.Range("AB2").Value = 5.129019193571377
$ws.Range("AC2").Value = 5.093556952227794
$ws.Range("AD2").Value = 5.063273301403831
$ws.Range("AE2").Value = 5.037476985802384
$ws.Range("AF2").Value = 5.019366639241145

$wb.Save()
